$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume/1h (E) columns for each coin row
$ws.Range("D2").Value = "67.935.63"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "3.255.34"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'583.70"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "'183.98"
$ws.Range("E6").Value = "  +4.14%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.598"
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("E9").Value = "  +3.74%  "
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("D12").Value = "3.824.86"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "'28.52"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").Value = "67.937.02"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("E16").Value = "  +2.66%  "
$ws.Range("D17").Value = "3.257.13"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "'5.86"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").Value = "'13.60"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").Value = "'382.18"
$ws.Range("E20").Value = "  +3.51%  "
$ws.Range("D21").Value = "'7.68"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'71.36"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "'9.83"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("E27").Value = "  +2.23%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").Value = "'5.69"
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("D31").Value = "'7.28"
$ws.Range("E31").Value = "  +7.53%  "
$ws.Range("D32").Value = "'22.92"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'1.26"
$ws.Range("E34").Value = "  +2.21%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  +3.77%  "
$ws.Range("D36").Value = "'162.46"
$ws.Range("E36").Value = "  -4.84%  "
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("D39").Value = "'26.62"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").Value = "'6.74"
$ws.Range("E40").Value = "  +4.83%  "
$ws.Range("D41").Value = "'4.60"
$ws.Range("E41").Value = "  +6.63%  "
$ws.Range("D42").Value = "'2.60"
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("D43").Value = "'41.39"
$ws.Range("E43").Value = "  +2.47%  "
$ws.Range("D44").Value = "'25.46"
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("D45").Value = "'0.0687"
$ws.Range("E45").Value = "  +2.12%  "
$ws.Range("E48").Value = "  +2.51%  "
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").Value = "'0.995"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").Value = "'31.27"
$ws.Range("E51").Value = "  +3.17%  "

# Rows 46 and 47 swapped coin (Maker <-> Bittensor) and got new price/volume data
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'344.82"
$ws.Range("E46").Value = "  +0.71%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.640.16"
$ws.Range("E47").Value = "  -3.92%  "

Write-Host "Done applying cryptos update"
